# Fixed naive component forecaster bug - Presentation state 11.02.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused C2/C3 values (first two years no longer have a
# naive YoY forecast component in column C).
$ws.Range("C2").ClearContents()
$ws.Range("C3").ClearContents()

# Recomputed (corrected) naive forecaster values - tiny precision fixes
# versus the previous buggy outputs.
$ws.Range("C4").Value  = -0.8792832172735632
$ws.Range("C5").Value  = 0.9337833426867226
$ws.Range("E5").Value  = 2.776721259569026
$ws.Range("C6").Value  = 2.791140000794257
$ws.Range("E6").Value  = 1.397663935193183
$ws.Range("C7").Value  = 0.4451370000809529
$ws.Range("E8").Value  = 1.564494423159402
$ws.Range("C11").Value = 2.2044495746113
$ws.Range("E11").Value = 1.120909053655073
$ws.Range("E13").Value = 1.897633937626786
$ws.Range("E14").Value = 0.5018989042238076
$ws.Range("C15").Value = -3.258619210312896
$ws.Range("E15").Value = -1.055440300316746
$ws.Range("C16").Value = 0.4255262881966759
$ws.Range("E16").Value = 2.531417074021181
$ws.Range("E17").Value = 1.719850910752729
$ws.Range("C18").Value = -0.2814561130375703
$ws.Range("C19").Value = -0.6470065423293758
$ws.Range("E19").Value = 1.208720904184779
